# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The underlying database query that feeds this "Estado de Cuenta" sheet now
# returns the debtor's overdue periods (col E, "Periodo Mora") newest-first
# instead of oldest-first, so the whole data block (rows 17-64) is
# re-written in the new order. Row 16 (period 1610) is the first/most-recent
# record under the new sort and keeps its place; the remaining 48 periods
# (1610..2009, one per row from 17 to 64) are now listed from 2009 down to
# 1610. The "Valor Mora" date in column F travels along with its period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Col E ("Periodo Mora") for rows 17-64, newest period first.
$periodos = @(
    "2009", "2008", "2007", "2006", "2005", "2004", "2003", "2002", "2001",
    "1912", "1911", "1910", "1909", "1908", "1907", "1906", "1905", "1904",
    "1903", "1902", "1901",
    "1812", "1811", "1810", "1809", "1808", "1807", "1806", "1805", "1804",
    "1803", "1802", "1801",
    "1712", "1711", "1710", "1709", "1708", "1707", "1706", "1705", "1704",
    "1703", "1702", "1701",
    "1612", "1611", "1610"
)

# Col F ("Valor Mora" date serial) for rows 17-64, same new order.
$valores = @(
    38610, 41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368,
    41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368,
    41368, 41368, 41368,
    41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368,
    41368, 41368, 41368,
    41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368, 41368,
    41368, 41368, 41368,
    41368, 41368, 28957
)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
